# Weekly update: insert a new price-record row at the top of the data
# block (row 232), pushing all existing data rows down by one.
#
# This mirrors the source workbook's "Fruta / hortaliza, semanal" commit:
# a brand-new weekly observation (Fecha = 2022-08-22, serial 44795) gets
# prepended to the "Apio" series for Feria Lagunitas de Puerto Montt,
# and every later row shifts from row N to row N+1 (last row 306 -> 307).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 232:306 down to 233:307, carrying formatting with them.
$ws.Rows.Item(232).Insert()

# Populate the newly-inserted row 232 with the new weekly record.
$ws.Range("A232").Value = 4
$ws.Range("B232").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C232").Value = "Los Lagos"
$ws.Range("D232").Value = 44795
$ws.Range("E232").Value = 10
$ws.Range("F232").Value = 100112017
$ws.Range("G232").Value = "Apio"
$ws.Range("H232").Value = "Americana (o)"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 25
$ws.Range("K232").Value = 14000
$ws.Range("L232").Value = 14000
$ws.Range("M232").Value = 14000
$ws.Range("N232").Value = "$/docena de matas"
$ws.Range("O232").Value = "Región de Coquimbo"
$ws.Range("P232").Value = 2333
$ws.Range("Q232").Value = 6
$ws.Range("R232").Value = "Hortaliza"
